# Updates the "cryptos" price list on Sheet1: refresh Price (column D) and
# Volume(1h) (column E) values to the latest scraped snapshot.
# GitHub Actions scheduled update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numeric-looking price strings (e.g. "64.734.23",
# "1.00") that must stay text. Assigning such a string straight to .Value
# lets Excel's COM layer "helpfully" reinterpret it as a number (dropping
# the decimal formatting / the thousands dots). Force the range to Text,
# write the string, then drop the style back to Normal so the written cell
# ends up with no explicit style (matching the rest of the sheet) while
# keeping the literal text.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '64.734.23'
$ws.Range("E2").Value = '  -0.71%  '
Set-TextValue $ws.Range("D3") '3.507.54'
$ws.Range("E3").Value = '  -1.27%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws.Range("D5") '587.18'
$ws.Range("E5").Value = '  -1.85%  '
Set-TextValue $ws.Range("D6") '132.48'
$ws.Range("E6").Value = '  -0.65%  '
Set-TextValue $ws.Range("D7") '3.506.89'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("E10").Value = '  +0.82%  '
Set-TextValue $ws.Range("D11") '7.12'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("E12").Value = '  -0.39%  '
Set-TextValue $ws.Range("D13") '4.104.77'
$ws.Range("E13").Value = '  -1.20%  '
Set-TextValue $ws.Range("D14") '27.75'
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("E16").Value = '  +0.72%  '
Set-TextValue $ws.Range("D17") '3.507.94'
$ws.Range("E17").Value = '  -1.29%  '
Set-TextValue $ws.Range("D18") '64.778.36'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("E19").Value = '  +0.64%  '
Set-TextValue $ws.Range("D20") '14.26'
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("E21").Value = '  -2.30%  '
Set-TextValue $ws.Range("D22") '390.77'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  -0.08%  '
Set-TextValue $ws.Range("D24") '3.647.93'
$ws.Range("E24").Value = '  -1.31%  '
Set-TextValue $ws.Range("D25") '74.15'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -4.21%  '
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("E29").Value = '  -4.39%  '
$ws.Range("E30").Value = '  -0.04%  '
Set-TextValue $ws.Range("D31") '2.26'
$ws.Range("E31").Value = '  -0.57%  '
Set-TextValue $ws.Range("D32") '8.20'
$ws.Range("E32").Value = '  -4.12%  '
Set-TextValue $ws.Range("D33") '3.508.97'
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("E34").Value = '  +0.00%  '
Set-TextValue $ws.Range("D35") '23.94'
$ws.Range("E35").Value = '  -0.46%  '
Set-TextValue $ws.Range("D36") '0.147'
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("E37").Value = '  +3.67%  '
Set-TextValue $ws.Range("D38") '5.25'
Set-TextValue $ws.Range("D39") '171.79'
$ws.Range("E39").Value = '  +0.84%  '
Set-TextValue $ws.Range("D40") '6.96'
$ws.Range("E40").Value = '  +0.45%  '
Set-TextValue $ws.Range("D41") '0.0814'
$ws.Range("E41").Value = '  +0.47%  '
Set-TextValue $ws.Range("D43") '26.19'
$ws.Range("E43").Value = '  -1.17%  '
Set-TextValue $ws.Range("D44") '1.00'
$ws.Range("E44").Value = '  +0.04%  '
Set-TextValue $ws.Range("D45") '42.34'
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("E47").Value = '  -0.87%  '
$ws.Range("E48").Value = '  -0.42%  '
Set-TextValue $ws.Range("D49") '2.480.68'
$ws.Range("E49").Value = '  +0.87%  '
Set-TextValue $ws.Range("D50") '6.87'
$ws.Range("E50").Value = '  -0.53%  '
Set-TextValue $ws.Range("D51") '0.906'
$ws.Range("E51").Value = '  +4.14%  '
